$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.132.96"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.266.87"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.29"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.34"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.493"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.05"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  -6.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.64"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.57"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "2.619.00"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "2.262.88"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.781"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "42.094.20"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.68"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.11"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.95"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.67"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.77"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.14"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.47"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.59"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0684"
$ws.Range("E38").Value = "  -4.16%  "
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0983"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").Value = "1.965.50"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0275"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.45"
$ws.Range("E46").Value = "  -5.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.50"
$ws.Range("E47").Value = "  -6.93%  "
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("D49").Value = "2.492.57"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.10"
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.74"
$ws.Range("E51").Value = "  -3.62%  "
